$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "277.77"
Set-TextValue $ws.Range("E2") "6.39%"
Set-TextValue $ws.Range("G2") "7"
Set-TextValue $ws.Range("D3") "27.41"
Set-TextValue $ws.Range("E3") "1.66%"
Set-TextValue $ws.Range("G3") "7"
Set-TextValue $ws.Range("D4") "4.808"
Set-TextValue $ws.Range("E4") "2.16%"
Set-TextValue $ws.Range("G4") "7"
Set-TextValue $ws.Range("D5") "0.06244"
Set-TextValue $ws.Range("E5") "0.44%"
Set-TextValue $ws.Range("G5") "7"
Set-TextValue $ws.Range("D6") "6.905"
Set-TextValue $ws.Range("E6") "2.25%"
Set-TextValue $ws.Range("G6") "7"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
Set-TextValue $ws.Range("D7") "3.278"
Set-TextValue $ws.Range("E7") "3.17%"
Set-TextValue $ws.Range("G7") "7"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D8") "0.8807"
Set-TextValue $ws.Range("E8") "3.43%"
Set-TextValue $ws.Range("G8") "7"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
Set-TextValue $ws.Range("D9") "0.9426"
Set-TextValue $ws.Range("E9") "3.03%"
Set-TextValue $ws.Range("G9") "7"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue $ws.Range("D10") "0.1449"
Set-TextValue $ws.Range("E10") "3.43%"
Set-TextValue $ws.Range("G10") "7"
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue $ws.Range("D11") "0.05234"
Set-TextValue $ws.Range("E11") "5.97%"
Set-TextValue $ws.Range("G11") "7"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue $ws.Range("D12") "0.07285"
Set-TextValue $ws.Range("E12") "2.85%"
Set-TextValue $ws.Range("G12") "7"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue $ws.Range("D13") "0.03158"
Set-TextValue $ws.Range("E13") "1.84%"
Set-TextValue $ws.Range("G13") "7"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue $ws.Range("D14") "0.09061"
Set-TextValue $ws.Range("E14") "0.07%"
Set-TextValue $ws.Range("G14") "7"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue $ws.Range("D15") "0.001562"
Set-TextValue $ws.Range("E15") "1.97%"
Set-TextValue $ws.Range("G15") "7"
$ws.Range("B16").Value = "One"
$ws.Range("C16").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue $ws.Range("D16") "0.0006282"
Set-TextValue $ws.Range("E16") "1.85%"
Set-TextValue $ws.Range("G16") "7"
$ws.Range("B17").Value = "TigerCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
Set-TextValue $ws.Range("D17") "0.005743"
Set-TextValue $ws.Range("E17") "-5.47%"
Set-TextValue $ws.Range("G17") "7"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue $ws.Range("D18") "3.453"
Set-TextValue $ws.Range("E18") "0.33%"
Set-TextValue $ws.Range("G18") "7"
Set-TextValue $ws.Range("D19") "2.284"
Set-TextValue $ws.Range("E19") "6.50%"
Set-TextValue $ws.Range("G19") "7"
Set-TextValue $ws.Range("D20") "0.3095"
Set-TextValue $ws.Range("E20") "-0.38%"
Set-TextValue $ws.Range("G20") "7"
Set-TextValue $ws.Range("D21") "0.1294"
Set-TextValue $ws.Range("E21") "-1.21%"
Set-TextValue $ws.Range("G21") "7"
Set-TextValue $ws.Range("D22") "3.869"
Set-TextValue $ws.Range("E22") "-5.89%"
Set-TextValue $ws.Range("G22") "7"
Set-TextValue $ws.Range("E23") "2.48%"
Set-TextValue $ws.Range("G23") "7"
Set-TextValue $ws.Range("D24") "0.001176"
Set-TextValue $ws.Range("E24") "-2.34%"
Set-TextValue $ws.Range("G24") "7"
Set-TextValue $ws.Range("D25") "0.004268"
Set-TextValue $ws.Range("E25") "4.58%"
Set-TextValue $ws.Range("G25") "7"
Set-TextValue $ws.Range("D26") "0.0001201"
Set-TextValue $ws.Range("E26") "0.05%"
Set-TextValue $ws.Range("G26") "7"
Set-TextValue $ws.Range("G27") "7"
Set-TextValue $ws.Range("G28") "7"
Set-TextValue $ws.Range("G29") "7"
Set-TextValue $ws.Range("G30") "7"
Set-TextValue $ws.Range("G31") "7"
Set-TextValue $ws.Range("G32") "7"
Set-TextValue $ws.Range("G33") "7"
Set-TextValue $ws.Range("G34") "7"
Set-TextValue $ws.Range("G35") "7"
Set-TextValue $ws.Range("G36") "7"
Set-TextValue $ws.Range("G37") "7"
Set-TextValue $ws.Range("G38") "7"
Set-TextValue $ws.Range("G39") "7"
Set-TextValue $ws.Range("D40") "0.04030"
Set-TextValue $ws.Range("E40") "2.25%"
Set-TextValue $ws.Range("G40") "7"
Set-TextValue $ws.Range("D41") "0.006394"
Set-TextValue $ws.Range("E41") "54.62%"
Set-TextValue $ws.Range("G41") "7"
Set-TextValue $ws.Range("E42") "3.64%"
Set-TextValue $ws.Range("G42") "7"
Set-TextValue $ws.Range("D43") "0.002135"
Set-TextValue $ws.Range("E43") "-3.46%"
Set-TextValue $ws.Range("G43") "7"
Set-TextValue $ws.Range("D44") "0.01206"
Set-TextValue $ws.Range("E44") "-10.71%"
Set-TextValue $ws.Range("G44") "7"
Set-TextValue $ws.Range("D45") "0.00005088"
Set-TextValue $ws.Range("E45") "-1.45%"
Set-TextValue $ws.Range("G45") "7"
Set-TextValue $ws.Range("E46") "0.05%"
Set-TextValue $ws.Range("G46") "7"
Set-TextValue $ws.Range("D47") "2.377"
Set-TextValue $ws.Range("E47") "823.23%"
Set-TextValue $ws.Range("G47") "7"
Set-TextValue $ws.Range("G48") "7"
Set-TextValue $ws.Range("E49") "0.05%"
Set-TextValue $ws.Range("G49") "7"
Set-TextValue $ws.Range("E50") "0.05%"
Set-TextValue $ws.Range("G50") "7"
Set-TextValue $ws.Range("G51") "7"
